# Generate Report for Handback
# Adds a new handback record (f4099a79-e6bc-4717-ab7e-9f17458e2234.md) as row 4
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$guidFile   = "f4099a79-e6bc-4717-ab7e-9f17458e2234.md"
$pathFile   = "e2e\f4099a79-e6bc-4717-ab7e-9f17458e2234.md"
$ext        = ".md"
$statusSync = "Handed back: in sync with en-US"
$dateFmt    = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $guidFile
$wsOverview.Range("B4").Value = $pathFile
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $statusSync
$wsOverview.Range("F4").Value = $statusSync
$wsOverview.Range("G4").Value = "2016-08-25 08:48:13"
$wsOverview.Range("G4").NumberFormat = $dateFmt

$wsOverview.Range("B4").Style = "HyperLink"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4099a79e6bc4717ab7e9f17458e2234/e2e/f4099a79-e6bc-4717-ab7e-9f17458e2234.md", "", "", $pathFile)

$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $guidFile
$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $statusSync
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = "f4099a79-e6bc-4717-ab7e-9f17458e2234.fa697d251c1c7511f3329a28bc2a17066135f240.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-25 08:47:59"
$wsZh.Range("H4").NumberFormat = $dateFmt
$wsZh.Range("I4").Value = $guidFile
$wsZh.Range("J4").Value = "f4099a79-e6bc-4717-ab7e-9f17458e2234.fa697d251c1c7511f3329a28bc2a17066135f240.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-08-25 08:48:30"
$wsZh.Range("K4").NumberFormat = $dateFmt
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Range("I4").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4099a79e6bc4717ab7e9f17458e2234/e2e/f4099a79-e6bc-4717-ab7e-9f17458e2234.md", "", "", $guidFile)
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f4099a79e6bc4717ab7e9f17458e2234/e2e/f4099a79-e6bc-4717-ab7e-9f17458e2234.md", "", "", $guidFile)

$loZh = $wsZh.ListObjects.Item("zh-cn")
$loZh.Resize($wsZh.Range("A1:P4"))

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $guidFile
$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $statusSync
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = "f4099a79-e6bc-4717-ab7e-9f17458e2234.fa697d251c1c7511f3329a28bc2a17066135f240.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-25 08:48:13"
$wsDe.Range("H4").NumberFormat = $dateFmt
$wsDe.Range("I4").Value = $guidFile
$wsDe.Range("J4").Value = "f4099a79-e6bc-4717-ab7e-9f17458e2234.fa697d251c1c7511f3329a28bc2a17066135f240.de-de.xlf"
$wsDe.Range("K4").Value = "2016-08-25 08:48:38"
$wsDe.Range("K4").NumberFormat = $dateFmt
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Range("I4").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4099a79e6bc4717ab7e9f17458e2234/e2e/f4099a79-e6bc-4717-ab7e-9f17458e2234.md", "", "", $guidFile)
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f4099a79e6bc4717ab7e9f17458e2234/e2e/f4099a79-e6bc-4717-ab7e-9f17458e2234.md", "", "", $guidFile)

$loDe = $wsDe.ListObjects.Item("de-de")
$loDe.Resize($wsDe.Range("A1:P4"))
